# Auto-generated edit script applying the commit diff
# (Update automàtic: dades i banners [2026-02-23 07:50])
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = '2026-02-23 07:48:18'
$ws.Cells.Item(2, 14).Value = '0.5 °C 7:00 TU'
$ws.Cells.Item(2, 15).Value = '3.0 °C'
$ws.Cells.Item(3, 5).Value = '2026-02-23 07:48:20'
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value = '40%'
$ws.Cells.Item(3, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(4, 5).Value = '2026-02-23 07:48:22'
$ws.Cells.Item(4, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(4, 15).Value = '5.5 °C'
$ws.Cells.Item(5, 5).Value = '2026-02-23 07:48:25'
$ws.Cells.Item(5, 8).NumberFormat = "@"
$ws.Cells.Item(5, 8).Value = '34%'
$ws.Cells.Item(5, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(6, 5).Value = '2026-02-23 07:48:27'
$ws.Cells.Item(6, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(7, 5).Value = '2026-02-23 07:48:30'
$ws.Cells.Item(7, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(7, 15).Value = '11.8 °C'
$ws.Cells.Item(8, 5).Value = '2026-02-23 07:48:32'
$ws.Cells.Item(8, 11).Value = '0.2 MJ/m2'
$ws.Cells.Item(8, 15).Value = '12.7 °C'
$ws.Cells.Item(9, 5).Value = '2026-02-23 07:48:35'
$ws.Cells.Item(9, 8).NumberFormat = "@"
$ws.Cells.Item(9, 8).Value = '92%'
$ws.Cells.Item(9, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(9, 15).Value = '6.5 °C'
$ws.Cells.Item(10, 5).Value = '2026-02-23 07:48:37'
$ws.Cells.Item(10, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(10, 13).Value = '5.8 °C 7:27 TU'
$ws.Cells.Item(11, 5).Value = '2026-02-23 07:48:39'
$ws.Cells.Item(11, 15).Value = '2.5 °C'
$ws.Cells.Item(12, 5).Value = '2026-02-23 07:48:42'
$ws.Cells.Item(12, 15).Value = '5.1 °C'
$ws.Cells.Item(13, 5).Value = '2026-02-23 07:48:44'
$ws.Cells.Item(13, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(13, 12).Value = '11.5 km/h - 324º 7:08 TU'
$ws.Cells.Item(13, 15).Value = '-1.5 °C'
$ws.Cells.Item(14, 5).Value = '2026-02-23 07:48:47'
$ws.Cells.Item(14, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(14, 13).Value = '10.3 °C 7:29 TU'
$ws.Cells.Item(14, 15).Value = '8.8 °C'
$ws.Cells.Item(15, 5).Value = '2026-02-23 07:48:49'
$ws.Cells.Item(15, 15).Value = '6.3 °C'
$ws.Cells.Item(16, 5).Value = '2026-02-23 07:48:51'
$ws.Cells.Item(16, 8).NumberFormat = "@"
$ws.Cells.Item(16, 8).Value = '16%'
$ws.Cells.Item(16, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(16, 15).Value = '2.7 °C'
$ws.Cells.Item(17, 5).Value = '2026-02-23 07:48:53'
$ws.Cells.Item(17, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(17, 13).Value = '8.6 °C 7:29 TU'
$ws.Cells.Item(18, 5).Value = '2026-02-23 07:48:56'
$ws.Cells.Item(18, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(19, 5).Value = '2026-02-23 07:48:58'
$ws.Cells.Item(19, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(20, 5).Value = '2026-02-23 07:49:01'
$ws.Cells.Item(20, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(21, 5).Value = '2026-02-23 07:49:03'
$ws.Cells.Item(21, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(21, 15).Value = '3.1 °C'
$ws.Cells.Item(22, 5).Value = '2026-02-23 07:49:06'
$ws.Cells.Item(22, 8).NumberFormat = "@"
$ws.Cells.Item(22, 8).Value = '23%'
$ws.Cells.Item(22, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(22, 12).Value = '22.7 km/h - 327º 7:26 TU'
$ws.Cells.Item(23, 5).Value = '2026-02-23 07:49:08'
$ws.Cells.Item(23, 8).NumberFormat = "@"
$ws.Cells.Item(23, 8).Value = '23%'
$ws.Cells.Item(23, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(23, 15).Value = '2.1 °C'
$ws.Cells.Item(24, 5).Value = '2026-02-23 07:49:10'
$ws.Cells.Item(24, 10).Value = '1027.6 hPa'
$ws.Cells.Item(24, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(24, 15).Value = '1.9 °C'
$ws.Cells.Item(25, 5).Value = '2026-02-23 07:49:13'
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value = '28%'
$ws.Cells.Item(25, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(25, 13).Value = '4.9 °C 7:12 TU'
$ws.Cells.Item(26, 5).Value = '2026-02-23 07:49:15'
$ws.Cells.Item(26, 7).Value = '1 cm'
$ws.Cells.Item(26, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(27, 5).Value = '2026-02-23 07:49:18'
$ws.Cells.Item(27, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(28, 5).Value = '2026-02-23 07:49:20'
$ws.Cells.Item(28, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(29, 5).Value = '2026-02-23 07:49:23'
$ws.Cells.Item(29, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(29, 13).Value = '5.2 °C 7:29 TU'
$ws.Cells.Item(30, 5).Value = '2026-02-23 07:49:25'
$ws.Cells.Item(30, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(31, 5).Value = '2026-02-23 07:49:27'
$ws.Cells.Item(31, 10).Value = '1024.3 hPa'
$ws.Cells.Item(31, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(31, 14).Value = '13.4 °C 7:01 TU'
$ws.Cells.Item(31, 15).Value = '15.0 °C'
$ws.Cells.Item(32, 5).Value = '2026-02-23 07:49:30'
$ws.Cells.Item(32, 8).NumberFormat = "@"
$ws.Cells.Item(32, 8).Value = '95%'
$ws.Cells.Item(32, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(32, 12).Value = '9.0 km/h - 288º 7:26 TU'
$ws.Cells.Item(32, 13).Value = '6.0 °C 7:29 TU'
$ws.Cells.Item(32, 15).Value = '1.5 °C'
$ws.Cells.Item(33, 5).Value = '2026-02-23 07:49:32'
$ws.Cells.Item(33, 10).Value = '1029.8 hPa'
$ws.Cells.Item(33, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(33, 15).Value = '2.1 °C'
$ws.Cells.Item(34, 5).Value = '2026-02-23 07:49:35'
$ws.Cells.Item(34, 8).NumberFormat = "@"
$ws.Cells.Item(34, 8).Value = '44%'
$ws.Cells.Item(34, 12).Value = '24.8 km/h - 36º 7:29 TU'
$ws.Cells.Item(34, 13).Value = '5.6 °C 7:23 TU'
$ws.Cells.Item(34, 15).Value = '2.5 °C'
$ws.Cells.Item(35, 5).Value = '2026-02-23 07:49:37'
$ws.Cells.Item(35, 8).NumberFormat = "@"
$ws.Cells.Item(35, 8).Value = '42%'
$ws.Cells.Item(35, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(36, 5).Value = '2026-02-23 07:49:40'
$ws.Cells.Item(36, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(36, 12).Value = '6.5 km/h - 53º 7:07 TU'
$ws.Cells.Item(37, 5).Value = '2026-02-23 07:49:42'
$ws.Cells.Item(37, 8).NumberFormat = "@"
$ws.Cells.Item(37, 8).Value = '82%'
$ws.Cells.Item(37, 15).Value = '3.2 °C'
$ws.Cells.Item(38, 5).Value = '2026-02-23 07:49:45'
$ws.Cells.Item(38, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(39, 5).Value = '2026-02-23 07:49:47'
$ws.Cells.Item(39, 8).NumberFormat = "@"
$ws.Cells.Item(39, 8).Value = '22%'
$ws.Cells.Item(39, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(39, 12).Value = '36.0 km/h - 326º 7:22 TU'
$ws.Cells.Item(40, 5).Value = '2026-02-23 07:49:49'
$ws.Cells.Item(40, 15).Value = '1.6 °C'
$ws.Cells.Item(41, 5).Value = '2026-02-23 07:49:52'
$ws.Cells.Item(41, 8).NumberFormat = "@"
$ws.Cells.Item(41, 8).Value = '86%'
$ws.Cells.Item(41, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(42, 5).Value = '2026-02-23 07:49:54'
$ws.Cells.Item(42, 8).NumberFormat = "@"
$ws.Cells.Item(42, 8).Value = '98%'
$ws.Cells.Item(43, 5).Value = '2026-02-23 07:49:56'
$ws.Cells.Item(43, 11).Value = '0.1 MJ/m2'
$ws.Cells.Item(43, 15).Value = '3.4 °C'
$ws.Cells.Item(44, 5).Value = '2026-02-23 07:49:59'
$ws.Cells.Item(44, 8).NumberFormat = "@"
$ws.Cells.Item(44, 8).Value = '35%'
$ws.Cells.Item(45, 5).Value = '2026-02-23 07:50:01'
$ws.Cells.Item(45, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(46, 5).Value = '2026-02-23 07:50:04'
$ws.Cells.Item(46, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(46, 15).Value = '1.8 °C'
